$d = $word.ActiveDocument

# 1) The opening "DON'T FORGET..." paragraph becomes two new paragraphs:
#    "Fetch in app" followed by "Work in folder".
$d.Content.Find.Execute(
    "DON’T FORGET TO COMMIT AND PUSH!!!",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Fetch in app^pWork in folder",
    2
)

# 2) The bookmarked "test" paragraph now carries the "DON'T FORGET..." text
#    right after the (still hidden/untouched) bookmark, and "test" moves
#    into its own new paragraph directly below it.
$d.Content.Find.Execute(
    "test",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "DON’T FORGET TO COMMIT AND PUSH!!!^ptest",
    2
)

# 3) Split the moved "DON'T FORGET..." paragraph from the "test" paragraph
#    with a blank paragraph in between.
$d.Paragraphs(4).Range.InsertParagraphAfter()
